$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking Service IDs in column C must stay TEXT cells (as in the
# source workbook), so the leading apostrophe forces text entry; the style
# is then reset back to "Normal" so no stray number-format survives on the
# cell itself (matches the original, un-styled <c> elements).
$idCells = @("C2","C3","C4","C5","C11","C12","C13","C14","C24")
$idValues = @(
    "10293410",
    "10293415",
    "10293432",
    "10293446",
    "10293483",
    "10293457",
    "10293466",
    "10293476",
    "137200001"
)

for ($i = 0; $i -lt $idCells.Length; $i++) {
    $ws.Range($idCells[$i]).Value = "'" + $idValues[$i]
    $ws.Range($idCells[$i]).Style = "Normal"
}

# Fail-log text swap on row 3.
$ws.Range("F3").Value = 'Cannot invoke "org.openqa.selenium.WebElement.isDisplayed()" because "element" is null'
